$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.261.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.487.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.72%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.07"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +20.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.492.11"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.694"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.98%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +32.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.29"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +10.08%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.030.87"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.80"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.27"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.471.82"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.286.04"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.05"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.89"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000140"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +27.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.35"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.87"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.22"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "313.86"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.60"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.23"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.69"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.179"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.39"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "44.39"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +13.04%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.81"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0493"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.66"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.60"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.43%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.99"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "136.55"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.42"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.288"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.32"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.201.45"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.836.48"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.39%  "
